$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "Then install the newest version of the Microsoft JDBC Driver,
# as shown in the picture." -> "Then install the Microsoft JDBC Driver, as
# shown in the picture." plus a new bold/red sentence appended right after
# it (before the trailing space + hyperlink that were already there).
# ---------------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute( `
    "Then install the newest version of the Microsoft JDBC Driver, as shown in the picture.", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "Then install the Microsoft JDBC Driver, as shown in the picture.", 2)

# Locate the (now shortened) sentence again and collapse to its end so we
# can append the new runs right after it, and before the pre-existing
# " " + hyperlink that followed in the same paragraph.
$insPoint = $d.Content
$insPoint.Find.Execute("Then install the Microsoft JDBC Driver, as shown in the picture.")
$insPoint.Collapse(0)

# Plain space, same formatting as the sentence that precedes it.
$insPoint.InsertAfter(" ")
$insPoint.Collapse(0)

# New bold, red sentence.
$insPoint.InsertAfter("I included the correct version of the driver in this repo.")
$insPoint.Font.Bold = 1
$insPoint.Font.Color = 255
$insPoint.Collapse(0)

# The pre-existing " " run right after this sentence (the one that was
# already followed by the hyperlink) also becomes bold + red.
$trailingSpace = $d.Content
$trailingSpace.Find.Execute("I included the correct version of the driver in this repo.")
$trailingSpace.Collapse(0)
$trailingSpace.MoveEnd(1, 1)
$trailingSpace.Font.Bold = 1
$trailingSpace.Font.Color = 255

# ---------------------------------------------------------------------------
# Change 2: drop the spell-check proofErr markers that bracket
# "CreatePollAppDB.sql" (spellStart before it, spellEnd after it) while
# leaving the surrounding text/runs otherwise intact.
# ---------------------------------------------------------------------------

# Clear the trailing proofErr marker first: touch the last character of the
# word together with the first character right after it.
$endProbe = $d.Content
$endProbe.Find.Execute("CreatePollAppDB.sql ")
$endRng = $d.Range($endProbe.End - 2, $endProbe.End)
$endText = $endRng.Text
$endRng.Delete()
$endRng.InsertAfter($endText)

# Clear the leading proofErr marker: touch the last character right before
# the word together with its first character.
$startProbe = $d.Content
$startProbe.Find.Execute(" CreatePollAppDB.sql")
$startRng = $d.Range($startProbe.Start, $startProbe.Start + 2)
$startText = $startRng.Text
$startRng.Delete()
$startRng.InsertAfter($startText)
